$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "27.865.85"

$ws.Cells.Item(3, 4).Value = "1.629.25"
$ws.Cells.Item(3, 5).Value = "  -0.35%  "

Set-TextValue 4 4 "0.998"
$ws.Cells.Item(4, 5).Value = "  -0.35%  "

Set-TextValue 5 4 "211.34"
$ws.Cells.Item(5, 5).Value = "  -0.32%  "

Set-TextValue 6 4 "0.517"
$ws.Cells.Item(6, 5).Value = "  -1.28%  "

Set-TextValue 8 4 "23.30"
$ws.Cells.Item(8, 5).Value = "  -0.54%  "

Set-TextValue 9 4 "0.258"
$ws.Cells.Item(9, 5).Value = "  -0.25%  "

$ws.Cells.Item(10, 5).Value = "  -1.15%  "

$ws.Cells.Item(11, 5).Value = "  -0.22%  "

$ws.Cells.Item(12, 4).Value = "1.860.20"
$ws.Cells.Item(12, 5).Value = "  -0.34%  "

$ws.Cells.Item(13, 4).Value = "1.612.10"
$ws.Cells.Item(13, 5).Value = "  -1.12%  "

$ws.Cells.Item(14, 5).Value = "  -1.19%  "

Set-TextValue 15 4 "0.556"
$ws.Cells.Item(15, 5).Value = "  -1.45%  "

Set-TextValue 16 4 "64.96"
$ws.Cells.Item(16, 5).Value = "  -1.26%  "

$ws.Cells.Item(17, 4).Value = "27.883.77"
$ws.Cells.Item(17, 5).Value = "  -0.28%  "

Set-TextValue 18 4 "228.74"
$ws.Cells.Item(18, 5).Value = "  -1.08%  "

Set-TextValue 19 4 "7.62"
$ws.Cells.Item(19, 5).Value = "  +0.06%  "

$ws.Cells.Item(20, 5).Value = "  -1.06%  "

$ws.Cells.Item(21, 5).Value = "  -0.34%  "

$ws.Cells.Item(22, 5).Value = "  -0.43%  "

$ws.Cells.Item(23, 5).Value = "  -4.87%  "

$ws.Cells.Item(24, 5).Value = "  -0.66%  "

Set-TextValue 25 4 "155.43"
$ws.Cells.Item(25, 5).Value = "  +0.01%  "

$ws.Cells.Item(26, 5).Value = "  -0.37%  "

$ws.Cells.Item(27, 5).Value = "  -0.17%  "

Set-TextValue 28 4 "15.46"

Set-TextValue 29 4 "0.998"
$ws.Cells.Item(29, 5).Value = "  -0.35%  "

$ws.Cells.Item(30, 5).Value = "  -0.22%  "

$ws.Cells.Item(31, 5).Value = "  -0.30%  "

$ws.Cells.Item(32, 5).Value = "  +0.25%  "

$ws.Cells.Item(33, 4).Value = "1.417.63"
$ws.Cells.Item(33, 5).Value = "  +0.77%  "

$ws.Cells.Item(34, 5).Value = "  +0.69%  "

$ws.Cells.Item(35, 5).Value = "  +2.50%  "

Set-TextValue 36 4 "1.00"
$ws.Cells.Item(36, 5).Value = "  -3.97%  "

$ws.Cells.Item(38, 5).Value = "  -0.85%  "

$ws.Cells.Item(39, 5).Value = "  -0.48%  "

Set-TextValue 40 4 "0.854"
$ws.Cells.Item(40, 5).Value = "  -1.42%  "

Set-TextValue 41 4 "1.01"
$ws.Cells.Item(41, 5).Value = "  -1.70%  "

Set-TextValue 42 4 "65.79"
$ws.Cells.Item(42, 5).Value = "  -1.57%  "

$ws.Cells.Item(43, 5).Value = "  -0.46%  "

$ws.Cells.Item(44, 5).Value = "  -0.54%  "

$ws.Cells.Item(45, 4).Value = "1.769.57"
$ws.Cells.Item(45, 5).Value = "  -0.37%  "

Set-TextValue 46 4 "2.12"
$ws.Cells.Item(46, 5).Value = "  -3.80%  "

Set-TextValue 47 4 "88.72"
$ws.Cells.Item(47, 5).Value = "  +0.49%  "

$ws.Cells.Item(48, 5).Value = "  +1.31%  "

Set-TextValue 49 4 "0.0502"
$ws.Cells.Item(49, 5).Value = "  -0.47%  "

Set-TextValue 50 4 "7.62"
$ws.Cells.Item(50, 5).Value = "  +1.68%  "

Set-TextValue 51 4 "0.997"
$ws.Cells.Item(51, 5).Value = "  -0.39%  "

Write-Host "Updated cryptos list"
